$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C values for rows 10-17 ("Top Global 50" playlist rows)
$ws.Range("C10").Value = 0.2516977999999999
$ws.Range("C11").Value = 0.6940400000000001
$ws.Range("C12").Value = 0.6271400000000001
$ws.Range("C13").Value = 0.001648606
$ws.Range("C14").Value = 0.175328
$ws.Range("C15").Value = 0.6
$ws.Range("C16").Value = 0.08866599999999998
$ws.Range("C17").Value = 0.51018
